# Applies the edits described by the commit:
#  1. Add a period "." right after "...Trapaça" in the "Recursos:" title
#     paragraph, as its own run (matching the Segoe UI / #444444 styling
#     used throughout that paragraph).
#  2. Bump the section's header distance (w:pgMar/@w:header) from 708 to
#     1278 twips (35.4pt -> 63.9pt).
#
$d = $word.ActiveDocument

# --- 1) Insert "." after "Trapaça" ------------------------------------
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $ptext = $p.Range.Text
    if ($ptext -like "Recursos:*Trapaça*") {
        $pr = $p.Range
        $insertPoint = $d.Range($pr.End - 1, $pr.End - 1)
        $insertPoint.InsertAfter(".")
        $insertPoint.Font.Name = "Segoe UI"
        $insertPoint.Font.NameAscii = "Segoe UI"
        $insertPoint.Font.NameBi = "Segoe UI"
        $insertPoint.Font.Color = 4473924
        break
    }
}

# --- 2) Header distance 708 -> 1278 twips (35.4pt -> 63.9pt) ----------
$d.PageSetup.HeaderDistance = 63.9
